$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, re-using the existing header style (same
# style as the other header cells, e.g. G1) instead of defining a new one.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the "Save" column values for rows 2-20 (no special style, matching
# the plain numeric cells in columns B-G).
$saveValues = @(1, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
